# Auto-generated Excel COM-interop script
# Applies numeric corrections to the Leve profit-tracking sheets (ALC, ARM, BSM, CRP, CUL, LTW, WVR)
# as produced by the scheduled market-data refresh runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1176.8
$ws.Range("J17").Value = 1176.8
$ws.Range("L17").Value = 3530.4
$ws.Range("N17").Value = -3866.4
$ws.Range("H33").Value = 1352
$ws.Range("I33").Value = 252.3077
$ws.Range("K33").Value = 252.3077
$ws.Range("M33").Value = -23.30770000000001
$ws.Range("H55").Value = 149.75
$ws.Range("J55").Value = 176
$ws.Range("L55").Value = 176
$ws.Range("N55").Value = -604
$ws.Range("H64").Value = 86049.836
$ws.Range("I64").Value = 102660.8
$ws.Range("J64").Value = 2995
$ws.Range("K64").Value = 102660.8
$ws.Range("L64").Value = 2995
$ws.Range("M64").Value = -102412.8
$ws.Range("N64").Value = -3491
$ws.Range("H67").Value = 86049.836
$ws.Range("I67").Value = 102660.8
$ws.Range("J67").Value = 2995
$ws.Range("K67").Value = 102660.8
$ws.Range("L67").Value = 2995
$ws.Range("M67").Value = -101802.8
$ws.Range("N67").Value = -4711
$ws.Range("H108").Value = 39744
$ws.Range("J108").Value = 39744
$ws.Range("L108").Value = 39744
$ws.Range("N108").Value = -47424
$ws.Range("H127").Value = 1821.0526
$ws.Range("I127").Value = 550
$ws.Range("J127").Value = 2160
$ws.Range("K127").Value = 1650
$ws.Range("L127").Value = 6480
$ws.Range("M127").Value = 3310
$ws.Range("N127").Value = -16400
$ws.Range("H132").Value = 6255972.5
$ws.Range("I132").Value = 6950705.5
$ws.Range("J132").Value = 3375
$ws.Range("K132").Value = 20852116.5
$ws.Range("L132").Value = 10125
$ws.Range("M132").Value = -20849586.5
$ws.Range("N132").Value = -15185
$ws.Range("H139").Value = 39137
$ws.Range("J139").Value = 39137
$ws.Range("L139").Value = 39137
$ws.Range("N139").Value = -49417

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7922.0537
$ws.Range("I32").Value = 7302.0576
$ws.Range("J32").Value = 15982
$ws.Range("K32").Value = 7302.0576
$ws.Range("L32").Value = 15982
$ws.Range("M32").Value = -7015.0576
$ws.Range("N32").Value = -16556
$ws.Range("H61").Value = 1712.6957
$ws.Range("I61").Value = 1590.1428
$ws.Range("K61").Value = 1590.1428
$ws.Range("M61").Value = -1378.1428
$ws.Range("H132").Value = 1966.5962
$ws.Range("I132").Value = 1963.8948
$ws.Range("J132").Value = 1973.9286
$ws.Range("K132").Value = 5891.6844
$ws.Range("L132").Value = 5921.7858
$ws.Range("M132").Value = -3361.6844
$ws.Range("N132").Value = -10981.7858
$ws.Range("H136").Value = 1712.6957
$ws.Range("I136").Value = 1590.1428
$ws.Range("K136").Value = 4770.428400000001
$ws.Range("M136").Value = -2220.428400000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 159640.86
$ws.Range("I86").Value = 222261.2
$ws.Range("J86").Value = 3090
$ws.Range("K86").Value = 222261.2
$ws.Range("L86").Value = 3090
$ws.Range("M86").Value = -221138.2
$ws.Range("N86").Value = -5336
$ws.Range("H89").Value = 159640.86
$ws.Range("I89").Value = 222261.2
$ws.Range("J89").Value = 3090
$ws.Range("K89").Value = 1111306
$ws.Range("L89").Value = 15450
$ws.Range("M89").Value = -1105690
$ws.Range("N89").Value = -26682
$ws.Range("H107").Value = 100002050
$ws.Range("I107").Value = 200001020
$ws.Range("K107").Value = 200001020
$ws.Range("M107").Value = -199999100

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2444.5637
$ws.Range("J31").Value = 2776.5898
$ws.Range("L31").Value = 2776.5898
$ws.Range("N31").Value = -3366.5898
$ws.Range("H34").Value = 2444.5637
$ws.Range("J34").Value = 2776.5898
$ws.Range("L34").Value = 2776.5898
$ws.Range("N34").Value = -3180.5898
$ws.Range("H132").Value = 6974.647
$ws.Range("I132").Value = 8328.556
$ws.Range("J132").Value = 5451.5
$ws.Range("K132").Value = 24985.668
$ws.Range("L132").Value = 16354.5
$ws.Range("M132").Value = -22455.668
$ws.Range("N132").Value = -21414.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1179.5
$ws.Range("I5").Value = 1049.8235
$ws.Range("J5").Value = 1250.6129
$ws.Range("K5").Value = 3149.4705
$ws.Range("L5").Value = 3751.8387
$ws.Range("M5").Value = -3037.4705
$ws.Range("N5").Value = -3975.8387
$ws.Range("H31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").ClearContents()
$ws.Range("H131").Value = 759.2525000000001
$ws.Range("J131").Value = 804.9666999999999
$ws.Range("L131").Value = 2414.9001
$ws.Range("N131").Value = -12494.9001
$ws.Range("H135").Value = 1179.5
$ws.Range("I135").Value = 1049.8235
$ws.Range("J135").Value = 1250.6129
$ws.Range("K135").Value = 9448.4115
$ws.Range("L135").Value = 11255.5161
$ws.Range("M135").Value = -6913.4115
$ws.Range("N135").Value = -16325.5161
$ws.Range("H139").Value = 2045.92
$ws.Range("I139").Value = 1027.7858
$ws.Range("J139").Value = 3341.7273
$ws.Range("K139").Value = 3083.3574
$ws.Range("L139").Value = 10025.1819
$ws.Range("M139").Value = 2056.6426
$ws.Range("N139").Value = -20305.1819

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3451.1562
$ws.Range("I7").Value = 4104.8
$ws.Range("J7").Value = 2874.4119
$ws.Range("K7").Value = 4104.8
$ws.Range("L7").Value = 2874.4119
$ws.Range("M7").Value = -3992.8
$ws.Range("N7").Value = -3098.4119
$ws.Range("H22").Value = 2711.1538
$ws.Range("I22").Value = 9800
$ws.Range("J22").Value = 2120.4167
$ws.Range("K22").Value = 9800
$ws.Range("L22").Value = 2120.4167
$ws.Range("M22").Value = -9505
$ws.Range("N22").Value = -2710.4167
$ws.Range("H27").Value = 2711.1538
$ws.Range("I27").Value = 9800
$ws.Range("J27").Value = 2120.4167
$ws.Range("K27").Value = 9800
$ws.Range("L27").Value = 2120.4167
$ws.Range("M27").Value = -9693
$ws.Range("N27").Value = -2334.4167
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("N32").ClearContents()
$ws.Range("H40").Value = 68746.92999999999
$ws.Range("I40").Value = 126562.375
$ws.Range("J40").Value = 2672.1428
$ws.Range("K40").Value = 126562.375
$ws.Range("L40").Value = 2672.1428
$ws.Range("M40").Value = -126426.375
$ws.Range("N40").Value = -2944.1428
$ws.Range("H122").Value = 2779.8
$ws.Range("I122").Value = 2779.8
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 8339.400000000001
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -5889.400000000001
$ws.Range("N122").ClearContents()
$ws.Range("H126").Value = 3451.1562
$ws.Range("I126").Value = 4104.8
$ws.Range("J126").Value = 2874.4119
$ws.Range("K126").Value = 12314.4
$ws.Range("L126").Value = 8623.235700000001
$ws.Range("M126").Value = -9844.400000000001
$ws.Range("N126").Value = -13563.2357
$ws.Range("H132").Value = 12360.6
$ws.Range("I132").Value = 23501.5
$ws.Range("J132").Value = 4933.3335
$ws.Range("K132").Value = 70504.5
$ws.Range("L132").Value = 14800.0005
$ws.Range("M132").Value = -67974.5
$ws.Range("N132").Value = -19860.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").ClearContents()
$ws.Range("H136").Value = 1284.4419
$ws.Range("I136").Value = 443.45715
$ws.Range("J136").Value = 4963.75
$ws.Range("K136").Value = 1330.37145
$ws.Range("L136").Value = 14891.25
$ws.Range("M136").Value = 1219.62855
$ws.Range("N136").Value = -19991.25

Write-Host "Applied all Leve profit corrections."
